$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.038640994678112
$ws.Cells.Item(2, 4).Value = 1.046467931157318
$ws.Cells.Item(2, 5).Value = 1.037277884794242
$ws.Cells.Item(2, 6).Value = 1.05554291047988
$ws.Cells.Item(2, 9).Value = 1.044069507245682
$ws.Cells.Item(2, 10).Value = 1.043737331682886
$ws.Cells.Item(2, 11).Value = 1.04923328360805
$ws.Cells.Item(2, 12).Value = 1.040069226611012
$ws.Cells.Item(2, 13).Value = 1.058283080269713

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.039494528143958
$ws.Cells.Item(3, 4).Value = 1.04713486524246
$ws.Cells.Item(3, 5).Value = 1.038001261542631
$ws.Cells.Item(3, 6).Value = 1.056347585341505
$ws.Cells.Item(3, 9).Value = 1.044307145050531
$ws.Cells.Item(3, 10).Value = 1.044236222902162
$ws.Cells.Item(3, 11).Value = 1.049712242374992
$ws.Cells.Item(3, 12).Value = 1.040602606856677
$ws.Cells.Item(3, 13).Value = 1.058901239591475

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.040047410332944
$ws.Cells.Item(4, 4).Value = 1.047566897909925
$ws.Cells.Item(4, 5).Value = 1.038470216357446
$ws.Cells.Item(4, 6).Value = 1.0568690870278
$ws.Cells.Item(4, 9).Value = 1.044459910944364
$ws.Cells.Item(4, 10).Value = 1.044558964392123
$ws.Cells.Item(4, 11).Value = 1.050021958790645
$ws.Cells.Item(4, 12).Value = 1.040947950140412
$ws.Cells.Item(4, 13).Value = 1.059301396817266

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.04027998141337
$ws.Cells.Item(5, 4).Value = 1.047748638014806
$ws.Cells.Item(5, 5).Value = 1.038667574257091
$ws.Cells.Item(5, 6).Value = 1.057088521517427
$ws.Cells.Item(5, 9).Value = 1.044523893106046
$ws.Cells.Item(5, 10).Value = 1.044694625813416
$ws.Cells.Item(5, 11).Value = 1.050152113925173
$ws.Cells.Item(5, 12).Value = 1.04109318151333
$ws.Cells.Item(5, 13).Value = 1.059469661278389

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.040319039249919
$ws.Cells.Item(6, 4).Value = 1.047779159587002
$ws.Cells.Item(6, 5).Value = 1.038700723738194
$ws.Cells.Item(6, 6).Value = 1.05712537692387
$ws.Cells.Item(6, 9).Value = 1.044534621858596
$ws.Cells.Item(6, 10).Value = 1.044717402806023
$ws.Cells.Item(6, 11).Value = 1.050173964595775
$ws.Cells.Item(6, 12).Value = 1.041117569336995
$ws.Cells.Item(6, 13).Value = 1.059497915825775

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.040050517413311
$ws.Cells.Item(7, 4).Value = 1.047569325886048
$ws.Cells.Item(7, 5).Value = 1.038472852643214
$ws.Cells.Item(7, 6).Value = 1.056872018357657
$ws.Cells.Item(7, 9).Value = 1.044460766823638
$ws.Cells.Item(7, 10).Value = 1.044560777183319
$ws.Cells.Item(7, 11).Value = 1.050023698126635
$ws.Cells.Item(7, 12).Value = 1.040949890538076
$ws.Cells.Item(7, 13).Value = 1.059303645026461

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.038929327403362
$ws.Cells.Item(8, 4).Value = 1.046693223857837
$ws.Cells.Item(8, 5).Value = 1.037522169722769
$ws.Cells.Item(8, 6).Value = 1.055814682513331
$ws.Cells.Item(8, 9).Value = 1.044150024959647
$ws.Cells.Item(8, 10).Value = 1.043905948967054
$ws.Cells.Item(8, 11).Value = 1.049395190993966
$ws.Cells.Item(8, 12).Value = 1.040249440538164
$ws.Cells.Item(8, 13).Value = 1.058491954623183

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.036958228503736
$ws.Cells.Item(9, 4).Value = 1.045153186975931
$ws.Cells.Item(9, 5).Value = 1.035853776398904
$ws.Cells.Item(9, 6).Value = 1.053957910967526
$ws.Cells.Item(9, 9).Value = 1.043594822964747
$ws.Cells.Item(9, 10).Value = 1.042751543562724
$ws.Cells.Item(9, 11).Value = 1.048286192167221
$ws.Cells.Item(9, 12).Value = 1.039016831250581
$ws.Cells.Item(9, 13).Value = 1.05706299944425

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.035647337271498
$ws.Cells.Item(10, 4).Value = 1.044129132773701
$ws.Cells.Item(10, 5).Value = 1.034746213536051
$ws.Cells.Item(10, 6).Value = 1.052724472659745
$ws.Cells.Item(10, 9).Value = 1.043219603827966
$ws.Cells.Item(10, 10).Value = 1.041981669605893
$ws.Cells.Item(10, 11).Value = 1.047545936146494
$ws.Cells.Item(10, 12).Value = 1.038196297730837
$ws.Cells.Item(10, 13).Value = 1.056111361924703

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.035080479736297
$ws.Cells.Item(11, 4).Value = 1.043686353520094
$ws.Cells.Item(11, 5).Value = 1.034267763660728
$ws.Cells.Item(11, 6).Value = 1.052191450069494
$ws.Cells.Item(11, 9).Value = 1.043055934115908
$ws.Cells.Item(11, 10).Value = 1.041648257454574
$ws.Cells.Item(11, 11).Value = 1.04722519443128
$ws.Cells.Item(11, 12).Value = 1.037841300805957
$ws.Cells.Item(11, 13).Value = 1.055699549115976

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.034870040429722
$ws.Cells.Item(12, 4).Value = 1.043521983972394
$ws.Cells.Item(12, 5).Value = 1.034090217987621
$ws.Cells.Item(12, 6).Value = 1.05199362358909
$ws.Cells.Item(12, 9).Value = 1.042994960819237
$ws.Cells.Item(12, 10).Value = 1.041524406850412
$ws.Cells.Item(12, 11).Value = 1.047106027042031
$ws.Cells.Item(12, 12).Value = 1.037709485516042
$ws.Cells.Item(12, 13).Value = 1.055546623027902

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.034915175063434
$ws.Cells.Item(13, 4).Value = 1.043557237317762
$ws.Cells.Item(13, 5).Value = 1.034128294328537
$ws.Cells.Item(13, 6).Value = 1.052036050685858
$ws.Cells.Item(13, 9).Value = 1.043008047888995
$ws.Cells.Item(13, 10).Value = 1.041550973498481
$ws.Cells.Item(13, 11).Value = 1.047131590171572
$ws.Cells.Item(13, 12).Value = 1.037737758222551
$ws.Cells.Item(13, 13).Value = 1.055579424378841

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.035063082365265
$ws.Cells.Item(14, 4).Value = 1.043672764664745
$ws.Cells.Item(14, 5).Value = 1.034253084162497
$ws.Cells.Item(14, 6).Value = 1.05217509435817
$ws.Cells.Item(14, 9).Value = 1.043050897692633
$ws.Cells.Item(14, 10).Value = 1.041638020049943
$ws.Cells.Item(14, 11).Value = 1.047215344614682
$ws.Cells.Item(14, 12).Value = 1.037830403958989
$ws.Cells.Item(14, 13).Value = 1.055686907382352

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.035154228462925
$ws.Cells.Item(15, 4).Value = 1.043743957965931
$ws.Cells.Item(15, 5).Value = 1.034329994132416
$ws.Cells.Item(15, 6).Value = 1.052260785257348
$ws.Cells.Item(15, 9).Value = 1.043077275165539
$ws.Cells.Item(15, 10).Value = 1.04169165148184
$ws.Cells.Item(15, 11).Value = 1.047266944603561
$ws.Cells.Item(15, 12).Value = 1.037887492241122
$ws.Cells.Item(15, 13).Value = 1.055753136493127

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.035684974034667
$ws.Cells.Item(16, 4).Value = 1.044158532273405
$ws.Cells.Item(16, 5).Value = 1.034777990687134
$ws.Cells.Item(16, 6).Value = 1.052759870203379
$ws.Cells.Item(16, 9).Value = 1.043230440890042
$ws.Cells.Item(16, 10).Value = 1.042003796092851
$ws.Cells.Item(16, 11).Value = 1.047567218489994
$ws.Cells.Item(16, 12).Value = 1.038219864135357
$ws.Cells.Item(16, 13).Value = 1.05613869804028

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.036018103180594
$ws.Cells.Item(17, 4).Value = 1.044418757184731
$ws.Cells.Item(17, 5).Value = 1.03505931148521
$ws.Cells.Item(17, 6).Value = 1.053073219325814
$ws.Cells.Item(17, 9).Value = 1.043326197709186
$ws.Cells.Item(17, 10).Value = 1.042199583231983
$ws.Cells.Item(17, 11).Value = 1.0477555181457
$ws.Cells.Item(17, 12).Value = 1.038428433401537
$ws.Cells.Item(17, 13).Value = 1.056380619338975

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.036212485675019
$ws.Cells.Item(18, 4).Value = 1.044570603891076
$ws.Cells.Item(18, 5).Value = 1.035223510261544
$ws.Cells.Item(18, 6).Value = 1.053256093130611
$ws.Cells.Item(18, 9).Value = 1.043381935460404
$ws.Cells.Item(18, 10).Value = 1.042313777489123
$ws.Cells.Item(18, 11).Value = 1.047865330187814
$ws.Cells.Item(18, 12).Value = 1.038550117131248
$ws.Cells.Item(18, 13).Value = 1.056521752197331

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.036278777568811
$ws.Cells.Item(19, 4).Value = 1.04462239014202
$ws.Cells.Item(19, 5).Value = 1.03527951625328
$ws.Cells.Item(19, 6).Value = 1.053318465709677
$ws.Cells.Item(19, 9).Value = 1.043400920963964
$ws.Cells.Item(19, 10).Value = 1.04235271388129
$ws.Cells.Item(19, 11).Value = 1.047902769830743
$ws.Cells.Item(19, 12).Value = 1.038591612982489
$ws.Cells.Item(19, 13).Value = 1.056569878921084

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.0359823539207
$ws.Cells.Item(20, 4).Value = 1.044390831097134
$ws.Cells.Item(20, 5).Value = 1.035029117126621
$ws.Cells.Item(20, 6).Value = 1.053039589307913
$ws.Cells.Item(20, 9).Value = 1.043315935857998
$ws.Cells.Item(20, 10).Value = 1.042178577639949
$ws.Cells.Item(20, 11).Value = 1.047735317429927
$ws.Cells.Item(20, 12).Value = 1.038406052903733
$ws.Cells.Item(20, 13).Value = 1.056354660952591

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.03501952411378
$ws.Cells.Item(21, 4).Value = 1.043638742010379
$ws.Cells.Item(21, 5).Value = 1.034216331906436
$ws.Cells.Item(21, 6).Value = 1.052134144979045
$ws.Cells.Item(21, 9).Value = 1.043038284430387
$ws.Cells.Item(21, 10).Value = 1.041612387178933
$ws.Cells.Item(21, 11).Value = 1.047190681826268
$ws.Cells.Item(21, 12).Value = 1.037803120807531
$ws.Cells.Item(21, 13).Value = 1.055655255212932

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.034414831479086
$ws.Cells.Item(22, 4).Value = 1.04316644383893
$ws.Cells.Item(22, 5).Value = 1.033706297117698
$ws.Cells.Item(22, 6).Value = 1.051565793521597
$ws.Cells.Item(22, 9).Value = 1.042862678268857
$ws.Cells.Item(22, 10).Value = 1.041256363868722
$ws.Cells.Item(22, 11).Value = 1.046848076994051
$ws.Cells.Item(22, 12).Value = 1.037424302456547
$ws.Cells.Item(22, 13).Value = 1.055215740530317

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.0347353258142
$ws.Cells.Item(23, 4).Value = 1.043416763403466
$ws.Cells.Item(23, 5).Value = 1.03397658123651
$ws.Cells.Item(23, 6).Value = 1.051866997887931
$ws.Cells.Item(23, 9).Value = 1.042955868318053
$ws.Cells.Item(23, 10).Value = 1.041445101633601
$ws.Cells.Item(23, 11).Value = 1.047029714051479
$ws.Cells.Item(23, 12).Value = 1.037625095275227
$ws.Cells.Item(23, 13).Value = 1.055448713298895

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.035998507244283
$ws.Cells.Item(24, 4).Value = 1.044403449502964
$ws.Cells.Item(24, 5).Value = 1.035042760320655
$ws.Cells.Item(24, 6).Value = 1.053054784948757
$ws.Cells.Item(24, 9).Value = 1.043320573104055
$ws.Cells.Item(24, 10).Value = 1.04218806917863
$ws.Cells.Item(24, 11).Value = 1.047744445325692
$ws.Cells.Item(24, 12).Value = 1.038416165598501
$ws.Cells.Item(24, 13).Value = 1.056366390355275

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.037467252191233
$ws.Cells.Item(25, 4).Value = 1.045550866266182
$ws.Cells.Item(25, 5).Value = 1.036284275545517
$ws.Cells.Item(25, 6).Value = 1.054437162046645
$ws.Cells.Item(25, 9).Value = 1.043739255752404
$ws.Cells.Item(25, 10).Value = 1.043050038212056
$ws.Cells.Item(25, 11).Value = 1.048573063397141
$ws.Cells.Item(25, 12).Value = 1.0393352836661
$ws.Cells.Item(25, 13).Value = 1.057432249334855
